$d = $word.ActiveDocument
$xmlFragment = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="3759DFB1" w14:textId="77777777" w:rsidR="00AA5AF7" w:rsidRPr="00C64508" w:rsidRDefault="00AA5AF7" w:rsidP="00AA5AF7"><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>The pizza dough</w:t></w:r></w:p><w:p w14:paraId="0DE75834" w14:textId="77777777" w:rsidR="00AA5AF7" w:rsidRPr="00C64508" w:rsidRDefault="00AA5AF7" w:rsidP="00AA5AF7"><w:pPr><w:pStyle w:val="Heading4"/></w:pPr><w:r><w:t>Ingredients</w:t></w:r></w:p><w:p w14:paraId="7F2F80A1" w14:textId="77777777" w:rsidR="00AA5AF7" w:rsidRDefault="00AA5AF7" w:rsidP="00AA5AF7"><w:r><w:t>400 g of flour</w:t></w:r></w:p><w:p w14:paraId="5C66941A" w14:textId="77777777" w:rsidR="00AA5AF7" w:rsidRDefault="00AA5AF7" w:rsidP="00AA5AF7"><w:r><w:t>1 sachet of bakery yeast</w:t></w:r></w:p><w:p w14:paraId="664C8DC6" w14:textId="77777777" w:rsidR="00AA5AF7" w:rsidRDefault="00AA5AF7" w:rsidP="00AA5AF7"><w:r><w:t>4 tablespoons of olive oil</w:t></w:r></w:p><w:p w14:paraId="21E7C717" w14:textId="77777777" w:rsidR="00AA5AF7" w:rsidRDefault="00AA5AF7" w:rsidP="00AA5AF7"><w:r><w:t>water</w:t></w:r></w:p><w:p w14:paraId="53CD6EE2" w14:textId="77777777" w:rsidR="00AA5AF7" w:rsidRDefault="00AA5AF7" w:rsidP="00AA5AF7"><w:r><w:t>1/2 Small teaspoon of salt</w:t></w:r></w:p><w:p w14:paraId="0B23B06B" w14:textId="77777777" w:rsidR="00AA5AF7" w:rsidRPr="00C64508" w:rsidRDefault="00AA5AF7" w:rsidP="00AA5AF7"><w:pPr><w:pStyle w:val="Heading4"/></w:pPr><w:r><w:t>Preparation</w:t></w:r></w:p><w:p w14:paraId="7333A3F4" w14:textId="77777777" w:rsidR="00AA5AF7" w:rsidRDefault="00AA5AF7" w:rsidP="00AA5AF7"><w:r><w:t>In the Kenwood bowl mix the flour with the yeast.</w:t></w:r></w:p><w:p w14:paraId="1F25B27E" w14:textId="77777777" w:rsidR="00AA5AF7" w:rsidRDefault="00AA5AF7" w:rsidP="00AA5AF7"><w:r><w:t>Add olive oil and salt, then while continuing to knead, add water.</w:t></w:r></w:p><w:p w14:paraId="043D4027" w14:textId="77777777" w:rsidR="00AA5AF7" w:rsidRDefault="00AA5AF7" w:rsidP="00AA5AF7"><w:r><w:t>The dough has the right consistency when the ball stands out from the walls of the bowl.</w:t></w:r></w:p><w:p w14:paraId="3CD2CA5F" w14:textId="77777777" w:rsidR="00AA5AF7" w:rsidRDefault="00AA5AF7" w:rsidP="00AA5AF7"><w:r><w:t>If there is too much water, add a little flour.</w:t></w:r></w:p><w:p w14:paraId="7DAA6CCF" w14:textId="77777777" w:rsidR="00AA5AF7" w:rsidRDefault="00AA5AF7" w:rsidP="00AA5AF7"><w:r><w:t>Pick up the dough in a ball at the bottom of the bowl, cover with a clean cloth and let up 2 hours at lukewarm temperature (24 ° C).</w:t></w:r></w:p><w:p w14:paraId="0E1D3A49" w14:textId="77777777" w:rsidR="00AA5AF7" w:rsidRDefault="00AA5AF7" w:rsidP="00AA5AF7"><w:r><w:t>When the dough has lifted, spread it in a pizza pan and allow to raise another 1/2 hour (optional)</w:t></w:r></w:p><w:p w14:paraId="1251CA61" w14:textId="77777777" w:rsidR="00AA5AF7" w:rsidRPr="00C64508" w:rsidRDefault="00AA5AF7" w:rsidP="00AA5AF7"><w:pPr><w:pStyle w:val="Heading4"/></w:pPr><w:r><w:t>Cooking</w:t></w:r></w:p><w:p w14:paraId="66CDA621" w14:textId="77777777" w:rsidR="00AA5AF7" w:rsidRDefault="00AA5AF7" w:rsidP="00AA5AF7"><w:r><w:t>The pizza dough is cooked quickly (15 minutes) in very hot preheated oven (240 ° C).</w:t></w:r></w:p><w:p w14:paraId="6673D7E6" w14:textId="77777777" w:rsidR="00571C58" w:rsidRPr="00AA5AF7" w:rsidRDefault="00571C58" w:rsidP="00AA5AF7"><w:r/></w:p><w:sectPr w:rsidR="00571C58" w:rsidRPr="00AA5AF7" w:rsidSect="003B7EFD"><w:pgSz w:w="11906" w:h="16838"/><w:pgMar w:top="1134" w:right="1133" w:bottom="1417" w:left="1417" w:header="708" w:footer="708" w:gutter="0"/><w:cols w:space="708"/><w:docGrid w:linePitch="360"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Content.InsertXML($xmlFragment)
Write-Output "done"
